# Update dashboards - 2025-11-26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: GDPNOW date cell loses its "latest update" yellow highlight ---
$ws.Range("C7").Interior.Pattern = -4142

# --- Row 13: UI Initial Claims (ICSA) refresh ---
$ws.Range("N13").Value = 45978
$ws.Range("Q13").Value = 216000
$ws.Range("R13").Value = 222000
$ws.Range("S13").Value = 228000
$ws.Range("T13").Value = 229000
$ws.Range("U13").Value = 220000

# --- Row 14: UI Continuing Claims (CCSA) refresh ---
$ws.Range("N14").Value = 45971
$ws.Range("Q14").Value = 1960000
$ws.Range("R14").Value = 1953000
$ws.Range("S14").Value = 1946000
$ws.Range("T14").Value = 1964000
$ws.Range("U14").Value = 1957000

# --- Row 22: PPI-FD M/M (PPIFIS) refresh - gains the yellow highlight ---
$ws.Range("N22").Interior.ColorIndex = 6
$ws.Range("N22").Value = 45901
$ws.Range("Q22").Value = 0.003100806343593332
$ws.Range("R22").Value = -0.001357686982925266
$ws.Range("S22").Value = 0.0081178572632572
$ws.Range("T22").Value = 0.0007422151749265637
$ws.Range("U22").Value = 0.003514212586162468

# --- Row 23: PPI-FD Y/Y (PPIFIS) refresh - gains the yellow highlight ---
$ws.Range("N23").Interior.ColorIndex = 6
$ws.Range("N23").Value = 45901
$ws.Range("Q23").Value = 0.02734717954345914
$ws.Range("R23").Value = 0.02721518987341764
$ws.Range("S23").Value = 0.0320979643678082
$ws.Range("T23").Value = 0.02406976503324603
$ws.Range("U23").Value = 0.02738206647949823

# --- Row 29: 5yr,5yr Forward (T5YIFR) refresh ---
$ws.Range("N29").Value = 45986
$ws.Range("Q29").Value = 2.17
$ws.Range("R29").Value = 2.16
$ws.Range("S29").Value = ""
$ws.Range("T29").Value = ""
$ws.Range("U29").Value = 2.16

# --- Row 30: 10yr TIPS (T10YIE) refresh ---
$ws.Range("N30").Value = 45986
$ws.Range("Q30").Value = 2.22
$ws.Range("R30").Value = 2.23
$ws.Range("S30").Value = ""
$ws.Range("T30").Value = ""
$ws.Range("U30").Value = 2.24

# --- Row 46-51: Exports/Imports/Trade Balance date cells lose highlight ---
$ws.Range("C46").Interior.Pattern = -4142
$ws.Range("C47").Interior.Pattern = -4142
$ws.Range("C48").Interior.Pattern = -4142
$ws.Range("C49").Interior.Pattern = -4142
$ws.Range("C50").Interior.Pattern = -4142
$ws.Range("C51").Interior.Pattern = -4142

# --- Row 47: FFR (DFF) latest date refresh ---
$ws.Range("N47").Value = 45985

# --- Row 48: 2y UST (DGS2) refresh ---
$ws.Range("N48").Value = 45985
$ws.Range("Q48").Value = 3.46
$ws.Range("R48").Value = ""
$ws.Range("S48").Value = ""
$ws.Range("T48").Value = 3.51
$ws.Range("U48").Value = 3.55

# --- Row 49: 5y UST (DGS5) refresh ---
$ws.Range("N49").Value = 45985
$ws.Range("Q49").Value = 3.61
$ws.Range("R49").Value = ""
$ws.Range("S49").Value = ""
$ws.Range("T49").Value = 3.62
$ws.Range("U49").Value = 3.68

# --- Row 50: 10y UST (DGS10) refresh ---
$ws.Range("N50").Value = 45985
$ws.Range("Q50").Value = 4.04
$ws.Range("R50").Value = ""
$ws.Range("S50").Value = ""
$ws.Range("T50").Value = 4.06
$ws.Range("U50").Value = 4.1

# --- Row 52: BAA (DBAA) refresh ---
$ws.Range("N52").Value = 45985
$ws.Range("Q52").Value = 5.84
$ws.Range("R52").Value = ""
$ws.Range("S52").Value = ""
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.9
